# Add "Turma" (class) and "Freq" (attendance) columns to the grades sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column C ("PI"), shifting
# PI..Q3 two columns to the right (C:K -> E:M).
$ws.Range("C1:D1").EntireColumn.Insert()

# Headers for the newly inserted columns.
$ws.Range("C1").Value = "Turma"

# Turma (class) values, entered so that "A" is encountered before "B"
# in the workbook's shared string table.
$ws.Range("C3").Value = "A"
$ws.Range("C6").Value = "A"
$ws.Range("C7").Value = "A"
$ws.Range("C2").Value = "B"
$ws.Range("C4").Value = "B"
$ws.Range("C5").Value = "B"

# Freq (attendance) header and values.
$ws.Range("D1").Value = "Freq"
$ws.Range("D2").Value = 80
$ws.Range("D3").Value = 90
$ws.Range("D4").Value = 70
$ws.Range("D5").Value = 100
$ws.Range("D6").Value = 63
$ws.Range("D7").Value = 75

# Corrected score values for Joaquim (row 6), now in shifted columns.
$ws.Range("J6").Value = 9
$ws.Range("M6").Value = 7

# Match the final cell selection recorded in the sheet view.
$ws.Range("D8").Select() | Out-Null
